$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: shift Designation_Name and Reporting_Manager_Code left,
# trim leading whitespace, and add new AttendanceLimit header in D1
$ws.Range("B1").Value = "Designation_Name"
$ws.Range("C1").Value = "Reporting_Manager_Code"
$ws.Range("D1").Value = "AttendanceLimit"

# Update data row values
$ws.Range("A2").Value = 10782
$ws.Range("B2").Value = "Account Executive"
$ws.Range("C2").Value = 10781
$ws.Range("D2").Value = 2

# Set column D width to match the other columns' style (~17.45 characters)
$ws.Columns.Item(4).ColumnWidth = 16.6

# Update the selected cell to D2, matching the new active cell in the saved file
$ws.Range("D2").Select() | Out-Null
